$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Row 31: GUI / tweakpane link / description
$ws.Range("B31").Value = "https://cocopon.github.io/tweakpane/ "
$ws.Range("A31").Value = "GUI"

# Row 32: npm install instructions
$ws.Range("B32").Value = "npm i tweakpane"
$ws.Range("A32").Value = "To Install (npm)"

# Column C descriptions
$ws.Range("C31").Value = "For setting up input sliders & GUI input that can be customised"
$ws.Range("C32").Value = "Installs tweakpane using npm"

# Add the hyperlink for the Tweakpane URL
$ws.Hyperlinks.Add($ws.Range("B31"), "https://cocopon.github.io/tweakpane/", "", "", "https://cocopon.github.io/tweakpane/ ") | Out-Null
$ws.Range("B31").Style = "Hyperlink"

# Update the selected cell to C32, matching the saved view state
$ws.Range("C32").Select()
